# Apply the commit's edits:
# 1. Clear contents of F2:F9 (the "Numerek" sequential numbers 1..8)
# 2. Move the active cell selection to C19

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2:F9").ClearContents()

$ws.Range("C19").Select()
